# Daily attendance processing - 2025-10-28 21:43:35
# Normalizes the "Recorded By" (column G) text so that the automated
# "System" account (and admin@admin.com) is listed first, ahead of the
# human/service account that also touched the record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Exact whole-cell replacements observed for the "Recorded By" column.
$replacements = @{
    "backup@backdoor.com, System"         = "System, backup@backdoor.com"
    "backup@backdoor.com, System, system" = "System, backup@backdoor.com, system"
    "dnasr281@gmail.com, System"          = "System, dnasr281@gmail.com"
    "dnasr281@gmail.com, admin@admin.com" = "admin@admin.com, dnasr281@gmail.com"
}

$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$lastRow = $usedRange.Rows.Count + $firstRow - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -eq $null) { continue }

    $text = [string]$val
    if ($replacements.ContainsKey($text)) {
        $cell.Value2 = $replacements[$text]
    }
}
